# Apply the "password" -> "hash_password" column rename and reset the
# encrypted password hash for every user except id_usuario 1, 2 and 3
# (rows 3-5 in the sheet) to the same new hash value.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("usuario")

# New unified password hash to apply to users 4 through 22 (rows 6-24)
$newHash = "dbb8d24a3166d17cd46539f4dca51ac945174b0ed8491859f690ddb24fe6cb74fbbf3338766d273903b8748da32c5a9c645ff0a2013c7412edc38b64cf8f2ec0"

for ($row = 6; $row -le 24; $row++) {
    $ws.Range("C$row").Value = $newHash
}

# Rename the column header in row 2 from "password" to "hash_password"
$ws.Range("C2").Value = "hash_password"

# Keep the active selection on C3 as in the saved workbook
$ws.Range("C3").Select()
